$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@("D2", "44.061.74")
    ,@("E2", "  -0.93%  ")
    ,@("D3", "2.201.26")
    ,@("E3", "  -2.14%  ")
    ,@("E4", "  +0.10%  ")
    ,@("D5", "294.16")
    ,@("E5", "  -4.39%  ")
    ,@("D6", "88.97")
    ,@("E6", "  -6.35%  ")
    ,@("E7", "  +0.59%  ")
    ,@("E8", "  +0.00%  ")
    ,@("D9", "0.483")
    ,@("E9", "  -8.08%  ")
    ,@("D10", "32.64")
    ,@("E10", "  -6.46%  ")
    ,@("D11", "0.0775")
    ,@("E11", "  -4.53%  ")
    ,@("D12", "0.102")
    ,@("E12", "  -1.74%  ")
    ,@("D13", "6.75")
    ,@("E13", "  -6.11%  ")
    ,@("B14", "WrappedEther")
    ,@("C14", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth")
    ,@("D14", "2.315.32")
    ,@("E14", "  -4.32%  ")
    ,@("B15", "WrappedliquidstakedEther2.0")
    ,@("C15", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth")
    ,@("D15", "2.539.20")
    ,@("E15", "  -2.08%  ")
    ,@("D16", "13.10")
    ,@("E16", "  -4.02%  ")
    ,@("D17", "0.765")
    ,@("E17", "  -8.64%  ")
    ,@("D18", "43.825.61")
    ,@("E18", "  -0.86%  ")
    ,@("D19", "0.0₃0886")
    ,@("E19", "  -7.98%  ")
    ,@("D20", "5.84")
    ,@("E20", "  -8.56%  ")
    ,@("D21", "10.83")
    ,@("E21", "  -12.80%  ")
    ,@("D22", "63.24")
    ,@("E22", "  -3.93%  ")
    ,@("D23", "229.80")
    ,@("E23", "  -3.57%  ")
    ,@("D24", "2.76")
    ,@("E24", "  -6.90%  ")
    ,@("E25", "  +0.03%  ")
    ,@("D26", "1.83")
    ,@("E26", "  -8.49%  ")
    ,@("D27", "2.21")
    ,@("E27", "  +0.10%  ")
    ,@("D28", "35.80")
    ,@("E28", "  -8.28%  ")
    ,@("D29", "9.20")
    ,@("E29", "  -6.70%  ")
    ,@("D30", "18.96")
    ,@("E30", "  -5.51%  ")
    ,@("D31", "147.32")
    ,@("E31", "  -3.99%  ")
    ,@("D32", "5.32")
    ,@("E32", "  -10.95%  ")
    ,@("D33", "2.50")
    ,@("E33", "  -5.14%  ")
    ,@("D34", "0.0737")
    ,@("E34", "  -7.92%  ")
    ,@("E35", "  -2.65%  ")
    ,@("D36", "2.86")
    ,@("E36", "  -8.79%  ")
    ,@("D37", "0.101")
    ,@("E37", "  -8.17%  ")
    ,@("D38", "1.65")
    ,@("E38", "  -7.79%  ")
    ,@("D39", "13.43")
    ,@("E39", "  -9.54%  ")
    ,@("B40", "VeChain")
    ,@("C40", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet")
    ,@("D40", "0.0282")
    ,@("E40", "  -6.93%  ")
    ,@("B41", "NEARProtocol")
    ,@("C41", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near")
    ,@("D41", "3.08")
    ,@("E41", "  -11.86%  ")
    ,@("D42", "3.49")
    ,@("E42", "  -8.82%  ")
    ,@("E43", "  -0.11%  ")
    ,@("D44", "1.746.56")
    ,@("E44", "  +0.47%  ")
    ,@("D45", "1.63")
    ,@("E45", "  -0.30%  ")
    ,@("D46", "68.65")
    ,@("E46", "  -0.90%  ")
    ,@("D47", "73.88")
    ,@("E47", "  -9.40%  ")
    ,@("D48", "0.172")
    ,@("E48", "  -10.37%  ")
    ,@("D49", "91.93")
    ,@("E49", "  -8.11%  ")
    ,@("D50", "2.420.90")
    ,@("E50", "  -2.11%  ")
    ,@("D51", "7.45")
    ,@("E51", "  -9.07%  ")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $cellVal = $u[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $cellVal
}
